$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.222.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.788.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.573.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.181.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.458.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
